$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (header row 4, data rows 5-7, blank separator row 3) currently
# runs through column M (year 2021). Extend it one column to the right (N)
# for the new year 2022, reusing the formatting already used for the 2021
# column so the new cells pick up the same styles (borders, number format,
# alignment, etc.) instead of Excel's plain default.
$ws.Range("M3:M7").Copy()
$ws.Range("N3:N7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new 2022 values.
$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 98.8
$ws.Range("N6").Value = 98
$ws.Range("N7").Value = 96.9

# Match the selection left behind in the edited workbook.
$ws.Range("O4").Select()
